# Update the "想去人数" (want-to-go count) figures on the "展览" sheet
# and on the "全部类型" sheet (which mirrors the same data plus one
# extra updated row, F32).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 414
$ws1.Range("F12").Value = 328
$ws1.Range("F21").Value = 308
$ws1.Range("F28").Value = 212
$ws1.Range("F29").Value = 4035
$ws1.Range("F34").Value = 121
$ws1.Range("F36").Value = 282

# --- Sheet "全部类型" (sheetId 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 414
$ws4.Range("F12").Value = 328
$ws4.Range("F21").Value = 308
$ws4.Range("F28").Value = 212
$ws4.Range("F29").Value = 4035
$ws4.Range("F32").Value = 246
$ws4.Range("F34").Value = 121
$ws4.Range("F36").Value = 282

$wb.Save()
